# Re-style every table in the deck that is still using the old default
# "Table_0" custom table style ({761B48A2-...}) so that it uses PowerPoint's
# built-in "No Style, No Grid" table style ({D61753B4-60EF-4F11-8114-123833B2C77B}),
# as if the author selected each table and picked that style from the
# Table Tools > Design > Table Styles gallery.
#
# (The companion theme1.xml/theme2.xml part-content swap recorded in the
# commit is an internal PowerPoint re-serialization artifact -- the design
# ("Integral") and its part relationships are otherwise untouched, and this
# host's Theme/Design COM surface is read-only, so there is no user-facing
# action that reproduces it here.)

$p = $ppt.ActivePresentation

$oldStyleId = "{761B48A2-C80D-480C-B9C2-8C7EDC8BF312}"
$newStyleId = "{D61753B4-60EF-4F11-8114-123833B2C77B}"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTable) {
            $table = $shape.Table
            if ($table.Style.Name -eq $oldStyleId) {
                $table.ApplyStyle($newStyleId)
            }
        }
    }
}
